$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.179.07"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "1.609.95"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.03"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3782"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.80"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3524"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08085"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.196"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.96"
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.353"
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.209"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001203"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("D17").Value = "1.611.55"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.01"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06915"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.512"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.13"
$ws.Range("E22").Value = "  -4.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.29"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").Value = "23.168.40"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.507"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.004"
$ws.Range("E26").Value = "  -8.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.81"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.99"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D31").Value = "1.789.60"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.065"
$ws.Range("E32").Value = "  +10.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.479"
$ws.Range("E33").Value = "  -5.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.095"
$ws.Range("E34").Value = "  -9.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.47"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02696"
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08719"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2444"
$ws.Range("E38").Value = "  -4.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06923"
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.807"
$ws.Range("E40").Value = "  -5.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.320"
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6841"
$ws.Range("E42").Value = "  -4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.92"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.18"
$ws.Range("E44").Value = "  -7.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6275"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.937"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.236"
$ws.Range("E48").Value = "  -4.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07848"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.75"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.160"
$ws.Range("E51").Value = "  -4.13%  "
